$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 152-153; existing rows 152.. shift down to 154..
$ws.Rows("152:153").Insert()

# New record: Terminal La Palmera de La Serena, Coquimbo - Uva, Autumn Royal
$ws.Range("A152").Value = 8
$ws.Range("B152").Value = "Terminal La Palmera de La Serena"
$ws.Range("C152").Value = "Coquimbo"
$ws.Range("D152").Value = 45075
$ws.Range("E152").Value = 4
$ws.Range("F152").Value = "Fruta"
$ws.Range("G152").Value = 100109
$ws.Range("H152").Value = "Uva"
$ws.Range("I152").Value = 100109001
$ws.Range("J152").Value = "Uva"
$ws.Range("K152").Value = "Autumn Royal"
$ws.Range("L152").Value = "Primera"
$ws.Range("M152").Value = 500
$ws.Range("N152").Value = 13000
$ws.Range("O152").Value = 14000
$ws.Range("P152").Value = 13500
$ws.Range("Q152").Value = "`$/bandeja 18 kilos"
$ws.Range("R152").Value = "Provincia del Elquí"
$ws.Range("S152").Value = 750
$ws.Range("T152").Value = 18

# New record: Terminal La Palmera de La Serena, Coquimbo - Uva, Red Globe
$ws.Range("A153").Value = 8
$ws.Range("B153").Value = "Terminal La Palmera de La Serena"
$ws.Range("C153").Value = "Coquimbo"
$ws.Range("D153").Value = 45075
$ws.Range("E153").Value = 4
$ws.Range("F153").Value = "Fruta"
$ws.Range("G153").Value = 100109
$ws.Range("H153").Value = "Uva"
$ws.Range("I153").Value = 100109001
$ws.Range("J153").Value = "Uva"
$ws.Range("K153").Value = "Red Globe"
$ws.Range("L153").Value = "Primera"
$ws.Range("M153").Value = 400
$ws.Range("N153").Value = 12000
$ws.Range("O153").Value = 13000
$ws.Range("P153").Value = 12500
$ws.Range("Q153").Value = "`$/bandeja 18 kilos"
$ws.Range("R153").Value = "Provincia del Elquí"
$ws.Range("S153").Value = 694
$ws.Range("T153").Value = 18
